$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Restyle the existing "trade_date" column (I2:I246) so it matches the
#    "datetime" column (B) format instead of the plain date-only format.
$ws.Range("I2:I246").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# 2. Append four new data rows (247-250) with the same layout/styling as
#    the existing rows: column B uses the datetime format, column I keeps
#    the original date-only format.
$newRows = @(
    @{ Row=247; A=179.27; B=46049; D=183.09; E=178.92; F=182.8;  H=2423458; I=46049 },
    @{ Row=248; A=180.98; B=46050; D=181.65; E=179.5;  F=179.5;  H=2270451; I=46050 },
    @{ Row=249; A=177.92; B=46051; D=182.26; E=176.05; F=181.3;  H=4004234; I=46051 },
    @{ Row=250; A=180.76; B=46052; D=181.89; E=178.34; F=179.4;  H=4980527; I=46052 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Range("A$row").Value = $r.A

    $ws.Range("B$row").Value = $r.B
    $ws.Range("B$row").NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Range("C$row").Value = "NSE"
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = "ITCHOT"
    $ws.Range("H$row").Value = $r.H

    $ws.Range("I$row").Value = $r.I
    $ws.Range("I$row").NumberFormat = "YYYY-MM-DD"

    $ws.Range("J$row").Value = "INE379A01028"
    $ws.Range("K$row").Value = "ITC Hotels Ltd"
    $ws.Range("L$row").Value = "ITCHOT"
    $ws.Range("M$row").Value = "BREEZE"
}
